$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1124.0364
$ws.Range("I129").Value = 616.2222
$ws.Range("J129").Value = 1371.081
$ws.Range("K129").Value = 1848.6666
$ws.Range("L129").Value = 4113.242999999999
$ws.Range("M129").Value = 3151.3334
$ws.Range("N129").Value = -14113.243
$ws.Range("H132").Value = 1668.3611
$ws.Range("I132").Value = 1120.5555
$ws.Range("J132").Value = 3311.7778
$ws.Range("K132").Value = 3361.6665
$ws.Range("L132").Value = 9935.3334
$ws.Range("M132").Value = -831.6664999999998
$ws.Range("N132").Value = -14995.3334
$ws.Range("H135").Value = 4323.5
$ws.Range("I135").Value = 5643.316
$ws.Range("K135").Value = 50789.844
$ws.Range("M135").Value = -48254.844
$ws.Range("H137").Value = 1275.125
$ws.Range("I137").Value = 647.5
$ws.Range("J137").Value = 1484.3334
$ws.Range("K137").Value = 1942.5
$ws.Range("L137").Value = 4453.0002
$ws.Range("M137").Value = 607.5
$ws.Range("N137").Value = -9553.0002
$ws.Range("H138").Value = 4803.278
$ws.Range("I138").Value = 3959.5293
$ws.Range("J138").Value = 5558.2104
$ws.Range("K138").Value = 11878.5879
$ws.Range("L138").Value = 16674.6312
$ws.Range("M138").Value = -6738.5879
$ws.Range("N138").Value = -26954.6312
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10674.605
$ws.Range("I32").Value = 3201.9106
$ws.Range("J32").Value = 38572.668
$ws.Range("K32").Value = 3201.9106
$ws.Range("L32").Value = 38572.668
$ws.Range("M32").Value = -2914.9106
$ws.Range("N32").Value = -39146.668
$ws.Range("H43").Value = 19333
$ws.Range("J43").Value = 19333
$ws.Range("L43").Value = 19333
$ws.Range("N43").Value = -19959
$ws.Range("H135").Value = 27903.375
$ws.Range("J135").Value = 27903.375
$ws.Range("L135").Value = 27903.375
$ws.Range("N135").Value = -38043.375
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 39800
$ws.Range("J135").Value = 39800
$ws.Range("L135").Value = 39800
$ws.Range("N135").Value = -49940
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1797.8793
$ws.Range("I31").Value = 1634.6957
$ws.Range("J31").Value = 1905.1143
$ws.Range("K31").Value = 1634.6957
$ws.Range("L31").Value = 1905.1143
$ws.Range("M31").Value = -1339.6957
$ws.Range("N31").Value = -2495.1143
$ws.Range("H34").Value = 1797.8793
$ws.Range("I34").Value = 1634.6957
$ws.Range("J34").Value = 1905.1143
$ws.Range("K34").Value = 1634.6957
$ws.Range("L34").Value = 1905.1143
$ws.Range("M34").Value = -1432.6957
$ws.Range("N34").Value = -2309.1143
$ws.Range("H134").Value = 1121.7142
$ws.Range("I134").Value = 1078.2778
$ws.Range("K134").Value = 3234.8334
$ws.Range("M134").Value = -699.8334000000004
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1449.8518
$ws.Range("I68").Value = 1094.5116
$ws.Range("J68").Value = 1851.9474
$ws.Range("K68").Value = 3283.5348
$ws.Range("L68").Value = 5555.8422
$ws.Range("M68").Value = -2472.5348
$ws.Range("N68").Value = -7177.8422
$ws.Range("H71").Value = 1449.8518
$ws.Range("I71").Value = 1094.5116
$ws.Range("J71").Value = 1851.9474
$ws.Range("K71").Value = 9850.6044
$ws.Range("L71").Value = 16667.5266
$ws.Range("M71").Value = -5794.6044
$ws.Range("N71").Value = -24779.5266
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6149.2856
$ws.Range("I70").Value = 5412.5
$ws.Range("K70").Value = 5412.5
$ws.Range("M70").Value = -5142.5
$ws.Range("H73").Value = 6149.2856
$ws.Range("I73").Value = 5412.5
$ws.Range("K73").Value = 5412.5
$ws.Range("M73").Value = -4476.5
$ws.Range("H124").Value = 42780
$ws.Range("J124").Value = 42780
$ws.Range("L124").Value = 42780
$ws.Range("N124").Value = -52600
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 406
$ws.Range("I22").Value = 176.66667
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 176.66667
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = 118.33333
$ws.Range("N22").Value = -1340
$ws.Range("H27").Value = 406
$ws.Range("I27").Value = 176.66667
$ws.Range("J27").Value = 750
$ws.Range("K27").Value = 176.66667
$ws.Range("L27").Value = 750
$ws.Range("M27").Value = -69.66667000000001
$ws.Range("N27").Value = -964
$ws.Range("H55").Value = 276.57144
$ws.Range("I55").Value = 233
$ws.Range("J55").Value = 385.5
$ws.Range("K55").Value = 233
$ws.Range("L55").Value = 385.5
$ws.Range("M55").Value = -60
$ws.Range("N55").Value = -731.5
$ws.Range("H136").Value = 2076.5881
$ws.Range("I136").Value = 1391.409
$ws.Range("J136").Value = 3332.75
$ws.Range("K136").Value = 4174.227000000001
$ws.Range("L136").Value = 9998.25
$ws.Range("M136").Value = -1624.227000000001
$ws.Range("N136").Value = -15098.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 8377
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 8377
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -9157
$ws.Range("H45").Value = 7686.25
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 7686.25
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -8668.25
$ws.Range("H74").Value = 5210
$ws.Range("I74").Value = 3500
$ws.Range("J74").Value = 5780
$ws.Range("K74").Value = 3500
$ws.Range("L74").Value = 5780
$ws.Range("M74").Value = -2564
$ws.Range("N74").Value = -7652
$ws.Range("H77").Value = 5210
$ws.Range("I77").Value = 3500
$ws.Range("J77").Value = 5780
$ws.Range("K77").Value = 10500
$ws.Range("L77").Value = 17340
$ws.Range("M77").Value = -5820
$ws.Range("N77").Value = -26700
$ws.Range("H132").Value = 983.325
$ws.Range("I132").Value = 718.125
$ws.Range("K132").Value = 2154.375
$ws.Range("M132").Value = 375.625
$ws.Range("H136").Value = 3491.6597
$ws.Range("I136").Value = 1299.2727
$ws.Range("J136").Value = 5420.96
$ws.Range("K136").Value = 3897.8181
$ws.Range("L136").Value = 16262.88
$ws.Range("M136").Value = -1347.8181
$ws.Range("N136").Value = -21362.88
